$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N ("Late") to make room for the
# upcoming "Variable Instalments" data - shifts old N->O ("Late") and old
# P->Q ("Outstanding") one column to the right.
$ws.Columns("N").Insert()

# The newly inserted column keeps the formatting of the column to its left,
# but its width should be a fixed 10 characters (not auto/best-fit).
$ws.Columns("N").ColumnWidth = 9.166667

# Restore the sheet's active selection.
$ws.Range("T8").Select() | Out-Null
